$wb = $excel.ActiveWorkbook

# --- Sheet "Appendix1" ---
$ws1 = $wb.Worksheets.Item("Appendix1")

# Header row (columns E, F, G change text)
$ws1.Range("E1").Value = "Số lượng trích dẫn trên Scopus"
$ws1.Range("F1").Value = "Số lượng trích dẫn trên Google Scholar"
$ws1.Range("G1").Value = "Ghi chú"

# Row 2 - new author data + numbers
$ws1.Range("B2").Value = "Trần Thị Thúy Nguyên"
$ws1.Range("C2").Value = "He130020"
$ws1.Range("D2").Value = "FPTUHN"
$ws1.Range("E2").Value = 135
$ws1.Range("F2").Value = 123

# Row 3 - re-assert author columns (string table is not de-duplicated by the
# runtime) plus the numbers that actually changed
$ws1.Range("B3").Value = "Trần Thị Thúy Nguyên"
$ws1.Range("C3").Value = "He130020"
$ws1.Range("D3").Value = "FPTUHN"
$ws1.Range("E3").Value = 123
$ws1.Range("F3").Value = 321

# --- Sheet "Appendix2" ---
$ws2 = $wb.Worksheets.Item("Appendix2")

# Row 2 - re-assert author columns + amount change
$ws2.Range("B2").Value = "Trần Thị Thúy Nguyên"
$ws2.Range("C2").Value = "He130020"
$ws2.Range("D2").Value = "FPTUHN"
$ws2.Range("E2").Value = 2000000

# Row 3 - author data changes + amount change
$ws2.Range("B3").Value = "Trần Thị Thúy Nguyên"
$ws2.Range("C3").Value = "He130020"
$ws2.Range("D3").Value = "FPTUHN"
$ws2.Range("E3").Value = 1000000
